# Natmi following Dr Hou advice
# Update sheet1 data: recompute ligand/receptor specificity metrics and add the
# "M2" sending-cluster rows (target cluster x M2) that were missing previously.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array holds the full A:T row values, in sheet-row order starting at row 2.
$rows = @(
    @("M1", "F9", "Lrp1", "ECs", 3, 1, 0.4394276666666667, 1.318283, 0.6309762956049347, 0.6309762956049347, 3, 1, 6.744736333333333, 20.234209, 0.01049273071342527, 0.01049273071342527, 2.963823749238556, 26.674413743147, 0.006620664356337202, 0.006620664356337204),
    @("M1", "F9", "Lrp1", "FAPs", 3, 1, 0.4394276666666667, 1.318283, 0.6309762956049347, 0.6309762956049347, 3, 1, 182.5316823333334, 547.595047, 0.283963033503136, 0.2839630335031361, 80.20947126047791, 721.8852413443011, 0.1791739429685487, 0.1791739429685488),
    @("M1", "F9", "Lrp1", "M1", 3, 1, 0.4394276666666667, 1.318283, 0.6309762956049347, 0.6309762956049347, 3, 1, 249.1329143333333, 747.398743, 0.3875740210972192, 0.3875740210972192, 109.4758952353632, 985.283057118269, 0.2445500201046322, 0.2445500201046322),
    @("M1", "F9", "Lrp1", "M2", 3, 1, 0.4394276666666667, 1.318283, 0.6309762956049347, 0.6309762956049347, 3, 1, 186.9310863333334, 560.793259, 0.2908071500393791, 0.2908071500393791, 82.1426910949219, 739.2842198542971, 0.1834924182672759, 0.1834924182672759),
    @("M1", "F9", "Lrp1", "sCs", 3, 1, 0.4394276666666667, 1.318283, 0.6309762956049347, 0.6309762956049347, 3, 1, 17.46044133333334, 52.38132400000001, 0.02716306464684043, 0.02716306464684043, 7.672600994076891, 69.05340894669202, 0.01713924990814074, 0.01713924990814074),
    @("M2", "F9", "Lrp1", "ECs", 2, 0.6666666666666666, 0.2569973333333334, 0.770992, 0.3690237043950653, 0.3690237043950653, 3, 1, 6.744736333333333, 20.234209, 0.01049273071342527, 0.01049273071342527, 1.733379251703111, 15.600413265328, 0.00387206635708807, 0.00387206635708807),
    @("M2", "F9", "Lrp1", "FAPs", 2, 0.6666666666666666, 0.2569973333333334, 0.770992, 0.3690237043950653, 0.3690237043950653, 3, 1, 182.5316823333334, 547.595047, 0.283963033503136, 0.2839630335031361, 46.91015560851378, 422.191400476624, 0.1047890905345873, 0.1047890905345873),
    @("M2", "F9", "Lrp1", "M1", 2, 0.6666666666666666, 0.2569973333333334, 0.770992, 0.3690237043950653, 0.3690237043950653, 3, 1, 249.1329143333333, 747.398743, 0.3875740210972192, 0.3875740210972192, 64.02649462922845, 576.238451663056, 0.143024000992587, 0.143024000992587),
    @("M2", "F9", "Lrp1", "M2", 2, 0.6666666666666666, 0.2569973333333334, 0.770992, 0.3690237043950653, 0.3690237043950653, 3, 1, 186.9310863333334, 560.793259, 0.2908071500393791, 0.2908071500393791, 48.04079070476979, 432.3671163429281, 0.1073147317721032, 0.1073147317721032),
    @("M2", "F9", "Lrp1", "sCs", 2, 0.6666666666666666, 0.2569973333333334, 0.770992, 0.3690237043950653, 0.3690237043950653, 3, 1, 17.46044133333334, 52.38132400000001, 0.02716306464684043, 0.02716306464684043, 4.487286861489778, 40.385581753408, 0.01002381473869969, 0.01002381473869969)
)

$r = 2
foreach ($rowVals in $rows) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $rowVals[$i]
    }
    $r = $r + 1
}
